# Trade #71 closed at 2026-02-17 08:57:35 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.53   # Current Capital
$wsSummary.Range("B4").Value = 0.54      # Total P&L $
$wsSummary.Range("B5").Value = 0.15      # Total P&L %
$wsSummary.Range("B6").Value = 71        # Total Trades
$wsSummary.Range("B8").Value = 29        # Losing Trades
$wsSummary.Range("B9").Value = 40.85     # Win Rate %

# --- Sheet: Strategy Status (row 4 = MarketMaking) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.53
$wsStatus.Range("D4").Value = 71
$wsStatus.Range("E4").Value = 0.54
$wsStatus.Range("F4").Value = 0.53
$wsStatus.Range("G4").Value = 40.85

# --- New trade row (row 72) shared by "All Trades" and "MarketMaking" ---
function Add-TradeRow($ws) {
    $row = 72

    $ws.Cells.Item($row, 1).Value = 71

    # Force text format on date/time columns so Excel doesn't coerce the
    # literal strings into date/time serial numbers.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "08:57:29"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.91
    $ws.Cells.Item($row, 7).Value = 0.87
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = -4.3956
    $ws.Cells.Item($row, 10).Value = -0.04
    $ws.Cells.Item($row, 11).Value = 100.53
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
